$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.647.89"
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("D3").Value = "'1.763.64"
$ws.Range("E3").Value = "  -2.10%  "
$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "  +1.04%  "
$ws.Range("D5").Value = "'326.31"
$ws.Range("E5").Value = "  -3.86%  "
$ws.Range("D6").Value = "'1.009"
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("D7").Value = "'0.4285"
$ws.Range("E7").Value = "  -8.41%  "
$ws.Range("D8").Value = "'0.3633"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'45.24"
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("D10").Value = "'1.124"
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").Value = "'0.07395"
$ws.Range("E11").Value = "  -2.65%  "
$ws.Range("D12").Value = "'1.012"
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").Value = "'21.81"
$ws.Range("E13").Value = "  -3.04%  "
$ws.Range("D14").Value = "'6.153"
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").Value = "'7.294"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").Value = "'1.766.00"
$ws.Range("E16").Value = "  -1.68%  "
$ws.Range("D17").Value = "'0.00001066"
$ws.Range("E17").Value = "  -2.01%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value = "'83.21"
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "'0.06253"
$ws.Range("E19").Value = "  -6.70%  "
$ws.Range("D20").Value = "'1.008"
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("D21").Value = "'16.97"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").Value = "'6.124"
$ws.Range("E22").Value = "  -4.30%  "
$ws.Range("D23").Value = "'27.730.28"
$ws.Range("E23").Value = "  -1.75%  "
$ws.Range("D24").Value = "'11.31"
$ws.Range("E24").Value = "  -5.26%  "
$ws.Range("D25").Value = "'2.403"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.384"
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'20.13"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("D28").Value = "'151.76"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").Value = "'1.969.59"
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("D30").Value = "'1.296"
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("D31").Value = "'128.00"
$ws.Range("E31").Value = "  -4.15%  "
$ws.Range("D32").Value = "'3.967"
$ws.Range("E32").Value = "  -2.52%  "
$ws.Range("D33").Value = "'5.672"
$ws.Range("E33").Value = "  -4.00%  "
$ws.Range("D34").Value = "'0.09069"
$ws.Range("E34").Value = "  -4.91%  "
$ws.Range("D35").Value = "'12.44"
$ws.Range("E35").Value = "  +2.62%  "
$ws.Range("D36").Value = "'0.2179"
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("D37").Value = "'0.06184"
$ws.Range("E37").Value = "  -1.64%  "
$ws.Range("D38").Value = "'0.6529"
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("D39").Value = "'0.02276"
$ws.Range("E39").Value = "  -4.11%  "
$ws.Range("D40").Value = "'5.075"
$ws.Range("E40").Value = "  -2.38%  "
$ws.Range("D41").Value = "'1.187"
$ws.Range("E41").Value = "  -2.27%  "
$ws.Range("D42").Value = "'1.427"
$ws.Range("E42").Value = "  -3.48%  "
$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").Value = "'1.009"
$ws.Range("E43").Value = "  +0.99%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'7.948"
$ws.Range("E44").Value = "  -1.68%  "
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("D46").Value = "'3.783"
$ws.Range("E46").Value = "  -2.27%  "
$ws.Range("D47").Value = "'0.5951"
$ws.Range("E47").Value = "  -2.45%  "
$ws.Range("D48").Value = "'125.17"
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("D49").Value = "'1.965"
$ws.Range("E49").Value = "  -3.27%  "
$ws.Range("D50").Value = "'0.06906"
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("D51").Value = "'1.113"
$ws.Range("E51").Value = "  -4.80%  "
